{"js": "// Replace the 25 \"dividend\u00f7divisor=quotient, remainder\" answer strings in\n// the table with their updated values, in document order.\nconst replacements = [\n  [\"742\u00f75=148, 2\", \"274\u00f76=45, 4\"],\n  [\"540\u00f77=77, 1\", \"542\u00f73=180, 2\"],\n  [\"943\u00f77=134, 5\", \"723\u00f79=80, 3\"],\n  [\"891\u00f78=111, 3\", \"846\u00f72=423, 0\"],\n  [\"189\u00f72=94, 1\", \"438\u00f74=109, 2\"],\n  [\"720\u00f79=80, 0\", \"172\u00f79=19, 1\"],\n  [\"842\u00f79=93, 5\", \"416\u00f75=83, 1\"],\n  [\"330\u00f77=47, 1\", \"317\u00f76=52, 5\"],\n  [\"364\u00f72=182, 0\", \"805\u00f74=201, 1\"],\n  [\"126\u00f72=63, 0\", \"983\u00f76=163, 5\"],\n  [\"663\u00f72=331, 1\", \"538\u00f79=59, 7\"],\n  [\"827\u00f78=103, 3\", \"546\u00f79=60, 6\"],\n  [\"978\u00f74=244, 2\", \"837\u00f78=104, 5\"],\n  [\"984\u00f77=140, 4\", \"447\u00f72=223, 1\"],\n  [\"242\u00f76=40, 2\", \"170\u00f76=28, 2\"],\n  [\"568\u00f75=113, 3\", \"348\u00f75=69, 3\"],\n  [\"855\u00f77=122, 1\", \"652\u00f74=163, 0\"],\n  [\"317\u00f78=39, 5\", \"802\u00f75=160, 2\"],\n  [\"624\u00f75=124, 4\", \"860\u00f74=215, 0\"],\n  [\"601\u00f76=100, 1\", \"167\u00f79=18, 5\"],\n  [\"107\u00f79=11, 8\", \"378\u00f73=126, 0\"],\n  [\"237\u00f74=59, 1\", \"605\u00f74=151, 1\"],\n  [\"258\u00f79=28, 6\", \"867\u00f79=96, 3\"],\n  [\"748\u00f72=374, 0\", \"658\u00f78=82, 2\"],\n  [\"259\u00f73=86, 1\", \"633\u00f76=105, 3\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  // Replace the first (and expected only) occurrence of this exact string.\n  found.items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 \"dividend\u00f7divisor=quotient, remainder\" answer strings in\n# the table with their updated values, in document order.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"742\u00f75=148, 2\", \"274\u00f76=45, 4\"),\n    @(\"540\u00f77=77, 1\", \"542\u00f73=180, 2\"),\n    @(\"943\u00f77=134, 5\", \"723\u00f79=80, 3\"),\n    @(\"891\u00f78=111, 3\", \"846\u00f72=423, 0\"),\n    @(\"189\u00f72=94, 1\", \"438\u00f74=109, 2\"),\n    @(\"720\u00f79=80, 0\", \"172\u00f79=19, 1\"),\n    @(\"842\u00f79=93, 5\", \"416\u00f75=83, 1\"),\n    @(\"330\u00f77=47, 1\", \"317\u00f76=52, 5\"),\n    @(\"364\u00f72=182, 0\", \"805\u00f74=201, 1\"),\n    @(\"126\u00f72=63, 0\", \"983\u00f76=163, 5\"),\n    @(\"663\u00f72=331, 1\", \"538\u00f79=59, 7\"),\n    @(\"827\u00f78=103, 3\", \"546\u00f79=60, 6\"),\n    @(\"978\u00f74=244, 2\", \"837\u00f78=104, 5\"),\n    @(\"984\u00f77=140, 4\", \"447\u00f72=223, 1\"),\n    @(\"242\u00f76=40, 2\", \"170\u00f76=28, 2\"),\n    @(\"568\u00f75=113, 3\", \"348\u00f75=69, 3\"),\n    @(\"855\u00f77=122, 1\", \"652\u00f74=163, 0\"),\n    @(\"317\u00f78=39, 5\", \"802\u00f75=160, 2\"),\n    @(\"624\u00f75=124, 4\", \"860\u00f74=215, 0\"),\n    @(\"601\u00f76=100, 1\", \"167\u00f79=18, 5\"),\n    @(\"107\u00f79=11, 8\", \"378\u00f73=126, 0\"),\n    @(\"237\u00f74=59, 1\", \"605\u00f74=151, 1\"),\n    @(\"258\u00f79=28, 6\", \"867\u00f79=96, 3\"),\n    @(\"748\u00f72=374, 0\", \"658\u00f78=82, 2\"),\n    @(\"259\u00f73=86, 1\", \"633\u00f76=105, 3\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
